# MitsosBarton2006Ex313 - "alpha_zero" stationary generator re-run.
# Updates the numeric results of the lider/follower restriction tables,
# the modified point, and the bf/BF vectors to the new solution values.
# Text cells stay as plain strings; pure-numeric-looking strings are
# written through a text-formatted cell (then the format is cleared) so
# they remain text (matching the workbook's original string layout)
# instead of being auto-coerced into number cells.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Restricciones_del_lider ---
$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
$wsLider.Range("A2").Value = "1.0499999999999998 - x"
Set-TextValue $wsLider.Range("B2") "-2.05"
Set-TextValue $wsLider.Range("D2") "0.24"
$wsLider.Range("A3").Value = "-1.05 + x"
Set-TextValue $wsLider.Range("B3") "0.050000000000000044"
Set-TextValue $wsLider.Range("D3") "0.72"

# --- Restricciones_del_follower ---
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
$wsFollower.Range("A2").Value = "-2.85 + y"
Set-TextValue $wsFollower.Range("B2") "1.85"
Set-TextValue $wsFollower.Range("D2") "0.47"
Set-TextValue $wsFollower.Range("E2") "0"
Set-TextValue $wsFollower.Range("F2") "7.9"
$wsFollower.Range("A3").Value = "2.85 - y"
Set-TextValue $wsFollower.Range("B3") "-3.85"
Set-TextValue $wsFollower.Range("D3") "0.88"
Set-TextValue $wsFollower.Range("E3") "0"
Set-TextValue $wsFollower.Range("F3") "7.9"

# --- Punto_modificado ---
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto.Range("A2") "1.05"
Set-TextValue $wsPunto.Range("B2") "2.85"

# --- Vector_bf ---
# NOTE: sheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) lookups are case-insensitive, so both names would
# resolve to the same (first) sheet. Use the 1-based tab index instead.
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf.Range("A2") "-1.4248750000000001"

# --- Vector_BF ---
$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF.Range("A2") "-1.48"
Set-TextValue $wsBF.Range("A3") "1.0"
